# UI Overhaul: Add Sidebar with Categories, Regenerate All Content
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column B (prompt -> shifts to D).
# This preserves existing formatting/shifts formulas automatically.
$ws.Range("B:C").Insert()

# Header row (row 1) - new headers for inserted columns, matching the
# same header style already used by A1/old-B1 (now D1): bold font,
# thin box border, centered horizontal / top vertical alignment.
$ws.Range("B1").Value = "category"
$ws.Range("C1").Value = "banner_type"
$ws.Range("B1:C1").Font.Bold = $true
$ws.Range("B1:C1").HorizontalAlignment = -4108
$ws.Range("B1:C1").VerticalAlignment = -4160
$ws.Range("B1:C1").Borders.LineStyle = 1
# D1 already has header style and "prompt" text preserved from the insert.

# Data for category (B) and banner_type (C) columns, row by row.
$categories = @{
    2 = "AI 그림"
    3 = "유튜브"
    4 = "업무 효율"
    5 = "AI 그림"
    6 = "생산성"
    7 = "디자인"
    8 = "블로그"
    9 = "SNS"
    10 = "코딩"
    11 = "마케팅"
}

$bannerTypes = @{
    2 = "tech"
    3 = "book"
    4 = "general"
    5 = "tech"
    6 = "general"
    7 = "tech"
    8 = "book"
    9 = "general"
    10 = "tech"
    11 = "book"
}

for ($row = 2; $row -le 11; $row++) {
    $topic = $ws.Range("A$row").Value2
    $ws.Range("B$row").Value = $categories[$row]
    $ws.Range("C$row").Value = $bannerTypes[$row]
    $ws.Range("D$row").Value = "Write a blog post about " + $topic
}
